$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Wrong_Entity_NonEvent_as_Event"
$ws.Range("B2").Value = 80

$ws.Range("A3").Value = "Wrong_Entity_Event_as_NonEvent"
$ws.Range("B3").Value = 49

$ws.Range("A4").Value = "Correct"
$ws.Range("B4").Value = 27

$ws.Range("A7").Value = "Wrong_Tag_S_as_B"
